# Auto-generated edit script: applies the row 13-16 data changes described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Force text storage (avoid auto-conversion of date-like / numeric-like strings),
    # then reset style back to Normal so we do not leave stray formatting behind.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

function Set-NumCell($addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-BoolCell($addr, $val) {
    $ws.Range($addr).Value = $val
}

function Clear-Cell($addr) {
    $ws.Range($addr).Value = ""
}

# ---- Row 13: update in place ----
Set-NumCell "A13" 112231588
Set-NumCell "B13" 93320
Set-NumCell "E13" 2818
Set-TextCell "F13" "Stubbspretmossa"
Set-TextCell "G13" "Herzogiella seligeri"
Set-TextCell "H13" "(Brid.) Z.Iwats."
Clear-Cell "I13"
Clear-Cell "J13"
Set-NumCell "Q13" 478719
Set-NumCell "R13" 6556487
Clear-Cell "AI13"
Set-TextCell "AJ13" "tall"
Set-TextCell "AK13" "Pinus sylvestris"
Set-TextCell "AM13" "Liggande död trädstam, markontakt"
Set-TextCell "AO13" "Horizontal, dead with ground contact # murken grov låga # Pinus sylvestris"
Set-TextCell "AY13" "Länsstyrelsen i Örebro län, inventering"

# ---- Row 14: update in place ----
Set-NumCell "A14" 112269209
Set-NumCell "B14" 56446
Set-TextCell "D14" "NT"
Set-NumCell "E14" 100049
Set-TextCell "F14" "Spillkråka"
Set-TextCell "G14" "Dryocopus martius"
Set-TextCell "H14" "(Linnaeus, 1758)"
Set-TextCell "I14" "1"
Clear-Cell "J14"
Set-TextCell "M14" "lockläte, övriga läten"
Set-NumCell "Q14" 478539
Set-NumCell "R14" 6556219
Set-NumCell "S14" 100
Clear-Cell "AF14"
Clear-Cell "AH14"
Clear-Cell "AJ14"
Clear-Cell "AK14"
Clear-Cell "AM14"
Clear-Cell "AO14"

# ---- Row 15: update in place ----
Set-NumCell "B15" 56575
Set-TextCell "AY15" "Länsstyrelsen i Örebro län, inventering"

# ---- Row 16: new row (append) ----
Set-NumCell "A16" 112205187
Set-NumCell "B16" 90796
Set-TextCell "C16" "Ovaliderad"
Set-TextCell "D16" "LC"
Set-NumCell "E16" 4363
Set-TextCell "F16" "Zontaggsvamp"
Set-TextCell "G16" "Hydnellum concrescens"
Set-TextCell "H16" "(Pers.) Banker"
Set-TextCell "I16" "1"
Set-TextCell "J16" "mycel"
Clear-Cell "K16"
Clear-Cell "N16"
Set-TextCell "P16" "Baggetorp, Nrk"
Set-NumCell "Q16" 478586
Set-NumCell "R16" 6556137
Set-NumCell "S16" 10
Set-TextCell "T16" "Örebro"
Set-TextCell "U16" "Lekeberg"
Set-TextCell "V16" "Närke"
Set-TextCell "W16" "Kvistbro"
Set-TextCell "Y16" "2023-09-19"
Set-TextCell "AA16" "2023-09-19"
Set-BoolCell "AD16" $false
Set-BoolCell "AE16" $false
Clear-Cell "AF16"
Set-BoolCell "AG16" $false
Set-TextCell "AH16" "Blåbärsbarrskog"
Set-TextCell "AI16" "i yta bökad av vildsvin"
Clear-Cell "AT16"
Set-TextCell "AW16" "Michael Andersson"
Set-TextCell "AX16" "Michael Andersson"
Set-TextCell "AY16" "Länsstyrelsen i Örebro län, inventering"
